# Update "想去人数" (F column) values on both the "展览" and "全部类型"
# worksheets, which contain mirrored data for the same events.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 2705
    $ws.Range("F4").Value = 586
    $ws.Range("F5").Value = 89
    $ws.Range("F7").Value = 870
}
